# Applies the commit "Refined metadata to be additional tab":
#  1. Refreshes the "time_taken" (column F) timestamps on the "data" sheet
#     for rows 2-58 to reflect the re-run query time.
#  2. Adds a new "metadata" worksheet after "data" describing the panel
#     query that produced the data (name, id, version, version_created,
#     query time, and the API request URL used).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timeTakenValues = @(
    "2021-10-05 14:34:27.731657",
    "2021-10-05 14:34:27.731665",
    "2021-10-05 14:34:27.731668",
    "2021-10-05 14:34:27.731670",
    "2021-10-05 14:34:27.731673",
    "2021-10-05 14:34:27.731676",
    "2021-10-05 14:34:27.731678",
    "2021-10-05 14:34:27.731681",
    "2021-10-05 14:34:27.731684",
    "2021-10-05 14:34:27.731686",
    "2021-10-05 14:34:27.731689",
    "2021-10-05 14:34:27.731691",
    "2021-10-05 14:34:27.731694",
    "2021-10-05 14:34:27.731696",
    "2021-10-05 14:34:27.731699",
    "2021-10-05 14:34:27.731701",
    "2021-10-05 14:34:27.731704",
    "2021-10-05 14:34:27.731706",
    "2021-10-05 14:34:27.731709",
    "2021-10-05 14:34:27.731711",
    "2021-10-05 14:34:27.731714",
    "2021-10-05 14:34:27.731716",
    "2021-10-05 14:34:27.731719",
    "2021-10-05 14:34:27.731721",
    "2021-10-05 14:34:27.731724",
    "2021-10-05 14:34:27.731726",
    "2021-10-05 14:34:27.731729",
    "2021-10-05 14:34:27.731731",
    "2021-10-05 14:34:27.731734",
    "2021-10-05 14:34:27.731736",
    "2021-10-05 14:34:27.731739",
    "2021-10-05 14:34:27.731741",
    "2021-10-05 14:34:27.731744",
    "2021-10-05 14:34:27.731747",
    "2021-10-05 14:34:27.731749",
    "2021-10-05 14:34:27.731752",
    "2021-10-05 14:34:27.731754",
    "2021-10-05 14:34:27.731756",
    "2021-10-05 14:34:27.731759",
    "2021-10-05 14:34:27.731761",
    "2021-10-05 14:34:27.731764",
    "2021-10-05 14:34:27.731767",
    "2021-10-05 14:34:27.731769",
    "2021-10-05 14:34:27.731772",
    "2021-10-05 14:34:27.731774",
    "2021-10-05 14:34:27.731776",
    "2021-10-05 14:34:27.731779",
    "2021-10-05 14:34:27.731781",
    "2021-10-05 14:34:27.731784",
    "2021-10-05 14:34:27.731786",
    "2021-10-05 14:34:27.731789",
    "2021-10-05 14:34:27.731791",
    "2021-10-05 14:34:27.731794",
    "2021-10-05 14:34:27.731796",
    "2021-10-05 14:34:27.731799",
    "2021-10-05 14:34:27.731801",
    "2021-10-05 14:34:27.731804"
)

for ($i = 0; $i -lt $timeTakenValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $timeTakenValues[$i]
}

# Add the new "metadata" worksheet positioned after the existing "data" sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$meta.Name = "metadata"

# Header row (bold, bordered, centered - same look as the "data" sheet header).
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row.
$idCell = $meta.Range("A2")
$idCell.Value = 0
$idCell.Font.Bold = $true
$idCell.Borders.LineStyle = 1
$idCell.HorizontalAlignment = -4108
$idCell.VerticalAlignment = -4160

$meta.Range("B2").Value = "Lysosomal Storage Disorder"
$meta.Range("C2").Value = 181

$versionCell = $meta.Range("D2")
$versionCell.NumberFormat = "@"
$versionCell.Value = "1.0"
$versionCell.Style = "Normal"

$meta.Range("E2").Value = "2021-04-14T04:56:42.538492Z"
$meta.Range("F2").Value = "2021-10-05 14:34:27.728184"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/181/?format=json"

Write-Output "metadata sheet added and time_taken column refreshed"
